$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width: 16 -> 22 (character units) ---
$ws.Columns("B").ColumnWidth = 21.17

# --- Apply the existing bordered/centered format (style used by A2:C3) to the
#     whole new block A5:C10 first, so everything starts from that baseline. ---
$ws.Range("A2:C2").Copy()
$ws.Range("A5:C10").PasteSpecial(-4122)   # xlPasteFormats

# A4 needs the same bordered/centered format too, but with the Hyperlink font.
# Apply the named "Hyperlink" style (brings in the right font) and then layer
# the border + centered alignment on top of it.
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("A4").Borders.LineStyle = 1
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4108

# B4 / C4 use the plain bordered/centered format (same as the rest).
$ws.Range("B2:C2").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)    # xlPasteFormats

# C5 additionally needs a date/time number format layered on top of the
# bordered/centered format already applied above.
$ws.Range("C5").NumberFormat = "m/d/yy h:mm"

# --- Values ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "alexf@gmail.com"
$ws.Range("C4").Value = "2024-07-28 00:23:02"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "alexfarouz@gmail.com"
$ws.Range("C5").Value = "2024-07-28 00:36:49"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "a@b.com"
$ws.Range("C6").Value = "2024-07-28 00:37:13"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "alexf@3.org"
$ws.Range("C7").Value = "2024-07-28 00:39:17"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "alexf@3.org"
$ws.Range("C8").Value = "2024-07-28 00:39:20"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "alexf@3.org"
$ws.Range("C9").Value = "2024-07-28 00:40:36"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "hi@m.com"
$ws.Range("C10").Value = "2024-07-28 00:40:48"
